$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: the workbook was originally saved from a different folder on the
# author's network drive, which Excel records automatically in the
# (read-only, Excel-managed) x15ac:absPath bookkeeping entry. That value is
# stamped by Excel itself from the real save location and isn't a
# settable property on the Workbook/Application object model, so it is
# not something a COM/VBA script can change directly.

# Update the version-number-like cells F2:F4 so they are stored as text
# (with a leading apostrophe / quote-prefix for the ones that look numeric).
# Order matters for shared-string table allocation, so write F4, F3, F2.
$ws.Range("F4").Value = "'1.10"
$ws.Range("F3").Value = "1.1.1"
$ws.Range("F2").Value = "'1.1"

# Update the active selection to match the target workbook.
$ws.Range("F3").Select()
